$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking text (prices like "1.002" or "0.3880").
# Excel auto-converts such strings to real numbers on assignment (dropping
# formatting such as trailing zeros), so force the cell to Text first, then
# restore the default (unstyled) cell style once the literal text is stored.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.853.82"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.96%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.647.88"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.22%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.78%  "

$ws.Range("E6").Value = "  -0.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3880"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.87%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3806"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.20%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "51.19"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.95%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.330"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.19%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.003"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.10%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08412"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.93%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.68"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.01%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.954"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.90%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.974"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.79%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001309"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.86%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.651.66"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.09%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.68"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.94%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06958"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.99%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.42%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.904"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.37%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.002"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.01%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.57"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.51%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.856.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.99%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.450"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.40%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.900"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.99%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.82"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.00%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "152.40"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.05%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.389"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.63%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "136.70"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.17%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.657"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.72%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.484"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.09%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.832.13"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.85%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08085"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.82%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9884"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.01%  "

$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02889"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.99%  "

$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.603"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.17%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2658"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.75%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "10.59"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.62%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09081"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.22%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.7503"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.45%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.38"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.89%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.414"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.81%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.46"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.53%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6882"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.80%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.420"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.60%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.093"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.69%  "

$ws.Range("E48").Value = "  -0.06%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08255"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.00%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.87"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.76%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.208"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.39%  "
